$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New leave-earned entries for periods that were previously blank (rows 27-33) ---
$ws.Range("C27").Value = 1.25
$ws.Range("C28").Value = 1.25
$ws.Range("C29").Value = 1.25
$ws.Range("C30").Value = 1.25
$ws.Range("C31").Value = 1.25
$ws.Range("C32").Value = 1.25
$ws.Range("C33").Value = 1.25

# --- New leave record for period ending 45231 (row 34): SL(2-0-0), 2 days absence w/ pay, remark ---
$ws.Range("B34").Value = "SL(2-0-0)"
$ws.Range("C34").Value = 1.25
$ws.Range("H34").Value = 2
$ws.Range("K34").Value = "11/28,29/2023"

# --- Running BALANCE formulas (column E) for rows 11-34 ---
$ws.Range("E11").Formula = "=SUM(C11,E10)-D11"
$ws.Range("E12:E34").Formula = "=SUM(C12,E11)-D12"

# --- Running BALANCE (2nd set, column I) formulas for rows 11-34 ---
$ws.Range("I11").Formula = "=SUM(G11,I10)-H11"
$ws.Range("I12:I34").Formula = "=SUM(G12,I11)-H12"

$wb.Application.Calculate()
